# "Add files via upload" re-save of Bessel beam information.xlsx
#
# The underlying numbers are untouched; the author just clarified the
# units in two of the tab names and left the selection on the
# "Solution of Z" sheet sitting on a different cell than before.

$wb = $excel.ActiveWorkbook

# Clarify the measurement units carried by each tab's data.
$wb.Worksheets.Item("Height information").Name = "Height information (mm)"
$wb.Worksheets.Item("Total phase").Name = "Total phase (rad)"

# "Solution of Z" keeps its name, but the last active cell on that sheet
# moved from E21 to I18 before the file was saved.
$ws3 = $wb.Worksheets.Item("Solution of Z")
[void]$ws3.Activate()
[void]$ws3.Range("I18").Select()
